# Generate Report for Archive
#
# Two changes, matching the "Ready for handoff" -> "In Translation" status
# rename plus the narrower Status/language columns now that the longer
# "Ready for handoff" text no longer needs to fit:
#
#   1. Update the status text everywhere it appears:
#        Overview!E2, Overview!F2, Overview!E3, Overview!F3  (zh-cn / de-de columns)
#        zh-cn!C2,    zh-cn!C3                                (Status column)
#        de-de!C2,    de-de!C3                                (Status column)
#
#   2. Narrow the columns that held that status text:
#        Overview columns E & F (zh-cn / de-de)
#        zh-cn column C  (Status)
#        de-de column C  (Status)

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$newColumnWidth = 13.4101845877511

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

foreach ($addr in @("E2", "F2", "E3", "F3")) {
    $cell = $wsOverview.Range($addr)
    if ($cell.Value() -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

foreach ($addr in @("C2", "C3")) {
    $cell = $wsZhCn.Range($addr)
    if ($cell.Value() -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

foreach ($addr in @("C2", "C3")) {
    $cell = $wsDeDe.Range($addr)
    if ($cell.Value() -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
